# Update column F ("dSF") values for specific rows, as part of a data
# repull / mean-calculation refresh. Column E ("dS0") is left unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    6  = -1
    7  = 1
    12 = 3
    13 = 0
    21 = -2
    23 = 3
    26 = 1
    27 = 2
    32 = -2
    35 = 0
    36 = 0
    37 = 1
    38 = -2
    40 = -3
    43 = 0
    44 = 3
    50 = -1
    51 = 2
    53 = 3
    55 = -3
    58 = 0
    61 = 2
    67 = -2
    69 = 0
    71 = 0
    72 = 1
    73 = 0
    80 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
